$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-5 to the reduced set of trials (2 for letters/phonemic, 2 for categories/semantic)
$ws.Cells.Item(2,1).Value2 = "Words that start with A"
$ws.Cells.Item(2,2).Value2 = "Phonemic"
$ws.Cells.Item(2,3).Value2 = "j"

$ws.Cells.Item(3,1).Value2 = "Words that start with S"
$ws.Cells.Item(3,2).Value2 = "Phonemic"
$ws.Cells.Item(3,3).Value2 = "p"

$ws.Cells.Item(4,1).Value2 = "Animals"
$ws.Cells.Item(4,2).Value2 = "Semantic"
$ws.Cells.Item(4,3).Value2 = "j"

$ws.Cells.Item(5,1).Value2 = "Occupations"
$ws.Cells.Item(5,2).Value2 = "Semantic"
$ws.Cells.Item(5,3).Value2 = "p"

# Remove the now-unused trailing rows (6-10) so the sheet shrinks to A1:C5
$ws.Range("A6:C10").ClearContents() | Out-Null

# Update selected cell to match the new saved view state
$ws.Range("A9").Select() | Out-Null
